$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.076.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.832.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.63"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6284"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07518"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2924"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.16"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07724"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.834.33"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.998"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6686"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.61"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009344"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -8.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.992"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.090.47"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.078.75"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.61"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "223.76"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.119"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.81"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1398"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.507"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.96"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.497"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05728"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +9.93%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.156"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.55%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.062"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.204"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7484"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.847"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.138"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.674"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.762"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.219.70"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01782"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.540"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.37%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8918"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.16"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.981.42"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000124"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.63"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07656"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +9.75%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5083"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4067"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.051"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.33%  "
